$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force each target cell to Text format before assigning, so that
# numeric-looking strings (prices) and percent-looking strings stay
# literal text instead of being auto-converted by Excel to Number/Percentage.
$cellValues = @{
    'D2' = '329.47'
    'E2' = '0.37%'
    'D3' = '44.40'
    'E3' = '0.57%'
    'D4' = '5.494'
    'E4' = '-1.45%'
    'D5' = '0.08030'
    'E5' = '-0.49%'
    'D6' = '2.068'
    'E6' = '7.75%'
    'D7' = '0.9537'
    'E7' = '0.78%'
    'D8' = '2.543'
    'E8' = '-0.18%'
    'D9' = '0.1144'
    'E9' = '-2.59%'
    'D10' = '0.1901'
    'E10' = '3.34%'
    'D11' = '10.17'
    'E11' = '5.94%'
    'D12' = '0.09888'
    'E12' = '2.36%'
    'E13' = '11.27%'
    'D14' = '0.1063'
    'E14' = '-0.49%'
    'D15' = '0.001267'
    'E15' = '-1.28%'
    'D16' = '0.04093'
    'E16' = '-2.88%'
    'D17' = '0.005934'
    'E17' = '0.90%'
    'D18' = '3.384'
    'E18' = '-4.36%'
    'D19' = '4.404'
    'E19' = '2.57%'
    'D20' = '0.3393'
    'E20' = '-2.98%'
    'D21' = '0.1383'
    'E21' = '1.42%'
    'D22' = '0.2578'
    'E22' = '-2.72%'
    'D23' = '0.001304'
    'E23' = '4.59%'
    'D24' = '0.004358'
    'E24' = '-2.94%'
    'D25' = '0.0001201'
    'E25' = '-4.72%'
    'D26' = '0.0003746'
    'E26' = '-6.14%'
    'D38' = '0.02600'
    'E38' = '-2.41%'
    'D39' = '0.05833'
    'E39' = '5.80%'
    'D40' = '0.007582'
    'E40' = '-0.06%'
    'D41' = '0.1403'
    'E41' = '0.19%'
    'D42' = '0.007359'
    'E42' = '3.93%'
    'D43' = '0.002009'
    'E43' = '-0.37%'
    'D44' = '0.009072'
    'E44' = '8.31%'
    'D45' = '0.00007054'
    'E45' = '2.05%'
    'D46' = '0.00000000751'
    'E46' = '0.04%'
    'D47' = '0.0005803'
    'E47' = '-0.15%'
    'D48' = '0.003532'
    'E48' = '55.51%'
    'E49' = '-38.09%'
    'D50' = '0.00002102'
    'E50' = '0.04%'
    'D51' = '0.0002002'
    'E51' = '0.04%'
}

foreach ($addr in $cellValues.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $cellValues[$addr]
}

Write-Host "Updated" $cellValues.Count "cells"
